# Update the "want to go" headcount column (F) for rows 2-10 on the
# sheets that carry this data: "展览" (sheet1) and "全部类型" (sheet4),
# which mirrors the same rows.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F (row 5 is unchanged).
$updates = @{
    2  = 644
    3  = 489
    4  = 34
    6  = 45
    7  = 39
    8  = 1915
    9  = 4035
    10 = 89
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
